# Created experiment order generation script
# Re-generates the randomized task-order tables for each condition sheet
# and re-labels the sheet tabs to match the freshly generated run IDs.

$wb = $excel.ActiveWorkbook

# xlPasteFormats = -4122 ; used to clone the existing bold/border/center
# "index column" style (style index 1 in styles.xml) onto newly added rows
# without having to hand-roll Font/Border/Alignment properties.
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Sheet 1 (was GNG) -> becomes the TOL task-order sheet, grows from
# 4 data rows (A1:B5) to 6 data rows (A1:B7).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TOL_TO-16515889958994648"

$ws1.Cells.Item(2,1).Value = 0
$ws1.Cells.Item(2,2).Value = "MM_stims-16515889958682482.csv"
$ws1.Cells.Item(3,1).Value = 1
$ws1.Cells.Item(3,2).Value = "ZM_stims-16515889958526227.csv"
$ws1.Cells.Item(4,1).Value = 2
$ws1.Cells.Item(4,2).Value = "MM_stims-16515889958838403.csv"
$ws1.Cells.Item(5,1).Value = 3
$ws1.Cells.Item(5,2).Value = "ZM_stims-16515889958682482.csv"

$ws1.Cells.Item(2,1).Copy()
$ws1.Cells.Item(6,1).PasteSpecial($xlPasteFormats)
$ws1.Cells.Item(6,1).Value = 4
$ws1.Cells.Item(6,2).Value = "MM_stims-16515889958994648.csv"

$ws1.Cells.Item(2,1).Copy()
$ws1.Cells.Item(7,1).PasteSpecial($xlPasteFormats)
$ws1.Cells.Item(7,1).Value = 5
$ws1.Cells.Item(7,2).Value = "ZM_stims-16515889958838403.csv"

# ---------------------------------------------------------------------
# Sheet 2 (NB) -> stays NB, same 9-data-row shape (A1:B10), only the
# generated filenames change.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-165158899759016"

$ws2.Cells.Item(2,1).Value = 0
$ws2.Cells.Item(2,2).Value = "ZB-match_7-1651588996286277.csv"
$ws2.Cells.Item(3,1).Value = 1
$ws2.Cells.Item(3,2).Value = "TB-16515889974005206.csv"
$ws2.Cells.Item(4,1).Value = 2
$ws2.Cells.Item(4,2).Value = "OB-16515889971373289.csv"
$ws2.Cells.Item(5,1).Value = 3
$ws2.Cells.Item(5,2).Value = "TB-1651588997574535.csv"
$ws2.Cells.Item(6,1).Value = 4
$ws2.Cells.Item(6,2).Value = "ZB-match_4-16515889965392883.csv"
$ws2.Cells.Item(7,1).Value = 5
$ws2.Cells.Item(7,2).Value = "OB-16515889971216981.csv"
$ws2.Cells.Item(8,1).Value = 6
$ws2.Cells.Item(8,2).Value = "OB-16515889968226764.csv"
$ws2.Cells.Item(9,1).Value = 7
$ws2.Cells.Item(9,2).Value = "ZB-match_0-16515889961903963.csv"
$ws2.Cells.Item(10,1).Value = 8
$ws2.Cells.Item(10,2).Value = "TB-16515889971712918.csv"

# ---------------------------------------------------------------------
# Sheet 3 (was RS, "eyes open"/"eyes closed") -> becomes the vSAT
# task-order sheet, grows from 2 data rows (A1:B3) to 4 data rows
# (A1:B5).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "vSAT_TO-16515889976526306"

$ws3.Cells.Item(2,1).Value = 0
$ws3.Cells.Item(2,2).Value = "vSAT_stims-16515889976370084.csv"
$ws3.Cells.Item(3,1).Value = 1
$ws3.Cells.Item(3,2).Value = "SAT_stims-16515889976058474.csv"

$ws3.Cells.Item(2,1).Copy()
$ws3.Cells.Item(4,1).PasteSpecial($xlPasteFormats)
$ws3.Cells.Item(4,1).Value = 2
$ws3.Cells.Item(4,2).Value = "SAT_stims-165158899759016.csv"

$ws3.Cells.Item(2,1).Copy()
$ws3.Cells.Item(5,1).PasteSpecial($xlPasteFormats)
$ws3.Cells.Item(5,1).Value = 3
$ws3.Cells.Item(5,2).Value = "vSAT_stims-16515889976213837.csv"

# ---------------------------------------------------------------------
# Sheet 4 (was TOL, MM/ZM stims) -> becomes the RS sheet, shrinks back
# down to the "eyes open" / "eyes closed" pair (A1:B3).
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "RS_TO-16515889976526306"
$ws4.Range("A4:B7").Clear()

$ws4.Cells.Item(2,1).Value = 0
$ws4.Cells.Item(2,2).Value = "eyes open"
$ws4.Cells.Item(3,1).Value = 1
$ws4.Cells.Item(3,2).Value = "eyes closed"

# ---------------------------------------------------------------------
# Sheet 5 (was vSAT) -> becomes the GNG sheet, same 4-data-row shape
# (A1:B5), only the generated filenames change.
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "GNG_TO-16515889976838837"

$ws5.Cells.Item(2,1).Value = 0
$ws5.Cells.Item(2,2).Value = "go_stims-16515889976526306.csv"
$ws5.Cells.Item(3,1).Value = 1
$ws5.Cells.Item(3,2).Value = "GNG_stims-1651588997668255.csv"
$ws5.Cells.Item(4,1).Value = 2
$ws5.Cells.Item(4,2).Value = "go_stims-1651588997668255.csv"
$ws5.Cells.Item(5,1).Value = 3
$ws5.Cells.Item(5,2).Value = "GNG_stims-16515889976838837.csv"
